# Remove the "** This course should provide" paragraph entirely
# (it sat between "navigate R help files" and the
# "this-course-will-not-address" bookmark/heading).

$d = $word.ActiveDocument

$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd("`r", "`a") -eq "** This course should provide") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    # Delete the whole paragraph, including its paragraph mark, so the
    # preceding and following paragraphs merge together seamlessly.
    $target.Range.Delete()
}
